# Auto-generated Excel COM-interop script applying the betting-odds update
# for France Ligue 1 (commit: Atualizacao de bases das ligas, 04-04-2024 23:22)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 151
$ws.Range('B151').Value = 6847967
$ws.Range('F151').Value = 'Metz'
$ws.Range('G151').Value = 'Montpellier'
$ws.Range('H151').Value = 0
$ws.Range('I151').Value = 1
$ws.Range('K151').Value = 3.2
$ws.Range('L151').Value = 3.5
$ws.Range('M151').Value = 2.15
$ws.Range('N151').Value = 3.3
$ws.Range('O151').Value = 3.25
$ws.Range('P151').Value = 2.25
$ws.Range('Q151').Value = 0.25
$ws.Range('R151').Value = 1.85
$ws.Range('S151').Value = 2
$ws.Range('U151').Value = 1.875
$ws.Range('V151').Value = 1.975
$ws.Range('Y151').Value = 1.25
$ws.Range('AA151').Value = 1
$ws.Range('AB151').Value = -1
$ws.Range('AC151').Value = 0.9750000000000001

# Row 152
$ws.Range('B152').Value = 6847961
$ws.Range('F152').Value = 'Lorient'
$ws.Range('G152').Value = 'Strasbourg'
$ws.Range('H152').Value = 1
$ws.Range('I152').Value = 2
$ws.Range('J152').Value = 'A'
$ws.Range('K152').Value = 2.3
$ws.Range('L152').Value = 3.25
$ws.Range('M152').Value = 3.1
$ws.Range('N152').Value = 2.375
$ws.Range('O152').Value = 3.2
$ws.Range('P152').Value = 3.1
$ws.Range('Q152').Value = -0.25
$ws.Range('R152').Value = 2.075
$ws.Range('S152').Value = 1.725
$ws.Range('U152').Value = 1.95
$ws.Range('V152').Value = 1.9
$ws.Range('X152').Value = -1
$ws.Range('Y152').Value = 2.1
$ws.Range('Z152').Value = -1
$ws.Range('AA152').Value = 0.7250000000000001
$ws.Range('AB152').Value = 0.95
$ws.Range('AC152').Value = -1

# Row 153
$ws.Range('B153').Value = 6847959
$ws.Range('F153').Value = 'Toulouse'
$ws.Range('G153').Value = 'Rennes'
$ws.Range('I153').Value = 0
$ws.Range('J153').Value = 'D'
$ws.Range('N153').Value = 3.2
$ws.Range('O153').Value = 3.4
$ws.Range('R153').Value = 1.9
$ws.Range('S153').Value = 1.95
$ws.Range('U153').Value = 1.825
$ws.Range('V153').Value = 2.025
$ws.Range('X153').Value = 2.4
$ws.Range('Y153').Value = -1
$ws.Range('Z153').Value = 0.45
$ws.Range('AA153').Value = -0.5
$ws.Range('AC153').Value = 1.025

# Row 178
$ws.Range('B178').Value = 6847989
$ws.Range('F178').Value = 'Lorient'
$ws.Range('G178').Value = 'Le Havre'
$ws.Range('H178').Value = 3
$ws.Range('I178').Value = 3
$ws.Range('K178').Value = 2.5
$ws.Range('L178').Value = 3.25
$ws.Range('M178').Value = 2.8
$ws.Range('N178').Value = 3.1
$ws.Range('O178').Value = 3.1
$ws.Range('P178').Value = 2.5
$ws.Range('Q178').Value = 0.25
$ws.Range('R178').Value = 1.81
$ws.Range('S178').Value = 2.125
$ws.Range('T178').Value = 2
$ws.Range('U178').Value = 2
$ws.Range('V178').Value = 1.9
$ws.Range('X178').Value = 2.1
$ws.Range('Z178').Value = 0.405
$ws.Range('AA178').Value = -0.5
$ws.Range('AB178').Value = 1
$ws.Range('AC178').Value = -1

# Row 179
$ws.Range('B179').Value = 6847992
$ws.Range('F179').Value = 'Clermont Foot'
$ws.Range('G179').Value = 'Strasbourg'
$ws.Range('H179').Value = 1
$ws.Range('I179').Value = 1
$ws.Range('K179').Value = 2.8
$ws.Range('L179').Value = 3.6
$ws.Range('M179').Value = 2.3
$ws.Range('N179').Value = 3
$ws.Range('O179').Value = 3.3
$ws.Range('P179').Value = 2.375
$ws.Range('R179').Value = 1.8
$ws.Range('S179').Value = 2.05
$ws.Range('T179').Value = 2.25
$ws.Range('U179').Value = 1.85
$ws.Range('V179').Value = 2
$ws.Range('X179').Value = 2.3
$ws.Range('Z179').Value = 0.4
$ws.Range('AB179').Value = -0.5
$ws.Range('AC179').Value = 0.5

# Row 180
$ws.Range('B180').Value = 6847993
$ws.Range('F180').Value = 'Reims'
$ws.Range('G180').Value = 'Nantes'
$ws.Range('H180').Value = 0
$ws.Range('I180').Value = 0
$ws.Range('K180').Value = 2
$ws.Range('L180').Value = 3.5
$ws.Range('M180').Value = 3.6
$ws.Range('N180').Value = 1.95
$ws.Range('O180').Value = 3.4
$ws.Range('P180').Value = 4
$ws.Range('Q180').Value = -0.5
$ws.Range('R180').Value = 2
$ws.Range('S180').Value = 1.93
$ws.Range('T180').Value = 2.5
$ws.Range('U180').Value = 2.07
$ws.Range('V180').Value = 1.83
$ws.Range('X180').Value = 2.4
$ws.Range('Z180').Value = -1
$ws.Range('AA180').Value = 0.9299999999999999
$ws.Range('AB180').Value = -1
$ws.Range('AC180').Value = 0.8300000000000001

# Row 187
$ws.Range('B187').Value = 7728864
$ws.Range('F187').Value = 'Reims'
$ws.Range('G187').Value = 'Toulouse'
$ws.Range('H187').Value = 2
$ws.Range('I187').Value = 3
$ws.Range('K187').Value = 1.75
$ws.Range('L187').Value = 3.5
$ws.Range('M187').Value = 4.75
$ws.Range('N187').Value = 1.85
$ws.Range('O187').Value = 3.4
$ws.Range('P187').Value = 4.5
$ws.Range('Q187').Value = -0.75
$ws.Range('T187').Value = 2.5
$ws.Range('U187').Value = 1.975
$ws.Range('V187').Value = 1.875
$ws.Range('Y187').Value = 3.5
$ws.Range('AB187').Value = 0.9750000000000001

# Row 189
$ws.Range('B189').Value = 7728866
$ws.Range('F189').Value = 'Metz'
$ws.Range('G189').Value = 'Lorient'
$ws.Range('H189').Value = 1
$ws.Range('I189').Value = 2
$ws.Range('K189').Value = 2.25
$ws.Range('L189').Value = 3.2
$ws.Range('M189').Value = 3.3
$ws.Range('N189').Value = 2.05
$ws.Range('O189').Value = 3.2
$ws.Range('P189').Value = 4
$ws.Range('Q189').Value = -0.5
$ws.Range('T189').Value = 2
$ws.Range('U189').Value = 1.9
$ws.Range('V189').Value = 1.95
$ws.Range('Y189').Value = 3
$ws.Range('AB189').Value = 0.8999999999999999

# Row 254
$ws.Range('N254').Value = 2.15
$ws.Range('P254').Value = 3.5
$ws.Range('R254').Value = 1.86
$ws.Range('S254').Value = 2.04
$ws.Range('T254').Value = 2.25
$ws.Range('U254').Value = 1.93
$ws.Range('V254').Value = 1.97

# Row 255
$ws.Range('N255').Value = 1.55
$ws.Range('P255').Value = 6.5
$ws.Range('R255').Value = 1.95
$ws.Range('S255').Value = 1.95
$ws.Range('U255').Value = 2.03
$ws.Range('V255').Value = 1.87

# Row 256
$ws.Range('O256').Value = 8.5
$ws.Range('P256').Value = 13
$ws.Range('R256').Value = 2.03
$ws.Range('S256').Value = 1.87

# Row 257
$ws.Range('O257').Value = 4
$ws.Range('P257').Value = 7.5
$ws.Range('R257').Value = 1.93
$ws.Range('S257').Value = 1.97
$ws.Range('U257').Value = 1.93
$ws.Range('V257').Value = 1.97

# Row 258
$ws.Range('O258').Value = 3.8
$ws.Range('R258').Value = 2.04
$ws.Range('S258').Value = 1.86

# Row 259
$ws.Range('N259').Value = 2.5
$ws.Range('O259').Value = 3.2
$ws.Range('R259').Value = 1.79
$ws.Range('S259').Value = 2.11
$ws.Range('T259').Value = 2.5
$ws.Range('U259').Value = 2.07
$ws.Range('V259').Value = 1.83

# Row 260
$ws.Range('O260').Value = 3.6
$ws.Range('U260').Value = 1.9
$ws.Range('V260').Value = 2

# Row 261
$ws.Range('O261').Value = 4
$ws.Range('P261').Value = 4
$ws.Range('U261').Value = 1.95
$ws.Range('V261').Value = 1.95

# Row 262
$ws.Range('N262').Value = 3.3
$ws.Range('R262').Value = 1.9
$ws.Range('S262').Value = 2
$ws.Range('U262').Value = 1.86
$ws.Range('V262').Value = 2.04
